$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.Range("K14")
$rng.Value = 7
$rng.Borders.Item(7).LineStyle = 1
$rng.Borders.Item(7).Color = 10855845
Write-Host "done1"
